$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the header style from an existing header cell (H1) to the new
# header cells I1 and J1 so they match the existing bold/bordered look.
$ws.Range("H1").Copy()
$ws.Range("I1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("J1").PasteSpecial(-4122)  # xlPasteFormats

# Header cells for new columns I and J
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Data for I (I0) and J (IF) columns, rows 2-31
$data = @{
    2  = @(1, 3)
    3  = @(2, 7)
    4  = @(2, 6)
    5  = @(1, 3)
    6  = @(5, 6)
    7  = @(1, 6)
    8  = @(1, 5)
    9  = @(1, 3)
    10 = @(3, 7)
    11 = @(1, 5)
    12 = @(1, 6)
    13 = @(1, 6)
    14 = @(2, 6)
    15 = @(1, 5)
    16 = @(1, 5)
    17 = @(1, 4)
    18 = @(1, 3)
    19 = @(1, 3)
    20 = @(1, 6)
    21 = @(1, 6)
    22 = @(1, 6)
    23 = @(1, 4)
    24 = @(1, 6)
    25 = @(1, 2)
    26 = @(1, 6)
    27 = @(1, 5)
    28 = @(5, 7)
    29 = @(4, 6)
    30 = @(1, 3)
    31 = @(1, 2)
}

foreach ($row in 2..31) {
    $vals = $data[$row]
    $ws.Cells.Item($row, 9).Value = $vals[0]
    $ws.Cells.Item($row, 10).Value = $vals[1]
}
